{"js": "// Replace each two-digit multiplication expression in the document with\n// its new value. Every \"before\" text is unique within the document, so a\n// single exact, case-sensitive search per pair is sufficient.\nconst replacements = [\n  [\"13\u00d786=\", \"36\u00d764=\"],\n  [\"93\u00d750=\", \"26\u00d752=\"],\n  [\"38\u00d753=\", \"48\u00d744=\"],\n  [\"16\u00d738=\", \"95\u00d749=\"],\n  [\"36\u00d754=\", \"94\u00d751=\"],\n  [\"69\u00d716=\", \"27\u00d712=\"],\n  [\"42\u00d738=\", \"65\u00d722=\"],\n  [\"50\u00d773=\", \"99\u00d789=\"],\n  [\"91\u00d779=\", \"56\u00d721=\"],\n  [\"69\u00d780=\", \"54\u00d725=\"],\n  [\"48\u00d772=\", \"27\u00d784=\"],\n  [\"92\u00d759=\", \"59\u00d763=\"],\n  [\"75\u00d748=\", \"33\u00d773=\"],\n  [\"41\u00d751=\", \"14\u00d745=\"],\n  [\"45\u00d780=\", \"82\u00d767=\"],\n  [\"50\u00d748=\", \"18\u00d793=\"],\n  [\"43\u00d773=\", \"15\u00d769=\"],\n  [\"27\u00d745=\", \"68\u00d774=\"],\n  [\"42\u00d783=\", \"98\u00d725=\"],\n  [\"83\u00d785=\", \"46\u00d719=\"],\n  [\"24\u00d743=\", \"59\u00d779=\"],\n  [\"39\u00d794=\", \"25\u00d760=\"],\n  [\"27\u00d772=\", \"46\u00d712=\"],\n  [\"74\u00d721=\", \"55\u00d744=\"],\n  [\"81\u00d768=\", \"95\u00d742=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression in the document with\n# its new value. Every \"before\" text is unique within the document, so a\n# single exact Find/Replace per pair is sufficient.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"13\u00d786=\", \"36\u00d764=\"),\n  @(\"93\u00d750=\", \"26\u00d752=\"),\n  @(\"38\u00d753=\", \"48\u00d744=\"),\n  @(\"16\u00d738=\", \"95\u00d749=\"),\n  @(\"36\u00d754=\", \"94\u00d751=\"),\n  @(\"69\u00d716=\", \"27\u00d712=\"),\n  @(\"42\u00d738=\", \"65\u00d722=\"),\n  @(\"50\u00d773=\", \"99\u00d789=\"),\n  @(\"91\u00d779=\", \"56\u00d721=\"),\n  @(\"69\u00d780=\", \"54\u00d725=\"),\n  @(\"48\u00d772=\", \"27\u00d784=\"),\n  @(\"92\u00d759=\", \"59\u00d763=\"),\n  @(\"75\u00d748=\", \"33\u00d773=\"),\n  @(\"41\u00d751=\", \"14\u00d745=\"),\n  @(\"45\u00d780=\", \"82\u00d767=\"),\n  @(\"50\u00d748=\", \"18\u00d793=\"),\n  @(\"43\u00d773=\", \"15\u00d769=\"),\n  @(\"27\u00d745=\", \"68\u00d774=\"),\n  @(\"42\u00d783=\", \"98\u00d725=\"),\n  @(\"83\u00d785=\", \"46\u00d719=\"),\n  @(\"24\u00d743=\", \"59\u00d779=\"),\n  @(\"39\u00d794=\", \"25\u00d760=\"),\n  @(\"27\u00d772=\", \"46\u00d712=\"),\n  @(\"74\u00d721=\", \"55\u00d744=\"),\n  @(\"81\u00d768=\", \"95\u00d742=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $old\n    $rng.Find.Replacement.Text = $new\n    $rng.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
